# Applies the per-cell numeric updates recorded in the commit diff.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets its
# changed cells written directly; one cell (ARM!M31) is cleared and
# one cell (ARM!N97) is newly populated, matching the diff's add/remove.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3222.389
$ws.Range("I76").Value = 3083.5833
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 3083.5833
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -2768.5833
$ws.Range("N76").Value = -4130
$ws.Range("H79").Value = 3222.389
$ws.Range("I79").Value = 3083.5833
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 3083.5833
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -1991.5833
$ws.Range("N79").Value = -5684
$ws.Range("H118").Value = 697.5
$ws.Range("I118").Value = 495
$ws.Range("K118").Value = 1485
$ws.Range("M118").Value = 172
$ws.Range("H129").Value = 877.1702
$ws.Range("I129").Value = 664.6667
$ws.Range("J129").Value = 891.6591
$ws.Range("K129").Value = 1994.0001
$ws.Range("L129").Value = 2674.9773
$ws.Range("M129").Value = 3005.9999
$ws.Range("N129").Value = -12674.9773
$ws.Range("H137").Value = 1909.2858
$ws.Range("I137").Value = 1969.7368
$ws.Range("J137").Value = 1837.5
$ws.Range("K137").Value = 5909.2104
$ws.Range("L137").Value = 5512.5
$ws.Range("M137").Value = -3359.2104
$ws.Range("N137").Value = -10612.5
$ws.Range("H138").Value = 3010.0425
$ws.Range("I138").Value = 1686.875
$ws.Range("J138").Value = 3692.9678
$ws.Range("K138").Value = 5060.625
$ws.Range("L138").Value = 11078.9034
$ws.Range("M138").Value = 79.375
$ws.Range("N138").Value = -21358.9034

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -385
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H32").Value = 2285.739
$ws.Range("I32").Value = 1395.05
$ws.Range("K32").Value = 1395.05
$ws.Range("M32").Value = -1108.05
$ws.Range("H61").Value = 4045.84
$ws.Range("I61").Value = 3607.8667
$ws.Range("J61").Value = 4702.8
$ws.Range("K61").Value = 3607.8667
$ws.Range("L61").Value = 4702.8
$ws.Range("M61").Value = -3395.8667
$ws.Range("N61").Value = -5126.8
$ws.Range("H74").Value = 1068.826
$ws.Range("J74").Value = 1278.8
$ws.Range("L74").Value = 1278.8
$ws.Range("N74").Value = -3026.8
$ws.Range("H77").Value = 1068.826
$ws.Range("J77").Value = 1278.8
$ws.Range("L77").Value = 6394
$ws.Range("N77").Value = -15130
$ws.Range("H97").Value = 5199.8
$ws.Range("J97").Value = 4239.5
$ws.Range("L97").Value = 4239.5
$ws.Range("N97").Value = -5231.5
$ws.Range("H124").Value = 11268.8
$ws.Range("J124").Value = 11268.8
$ws.Range("L124").Value = 11268.8
$ws.Range("N124").Value = -21088.8
$ws.Range("H136").Value = 4045.84
$ws.Range("I136").Value = 3607.8667
$ws.Range("J136").Value = 4702.8
$ws.Range("K136").Value = 10823.6001
$ws.Range("L136").Value = 14108.4
$ws.Range("M136").Value = -8273.6001
$ws.Range("N136").Value = -19208.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5000100
$ws.Range("I7").Value = 5000075
$ws.Range("K7").Value = 5000075
$ws.Range("M7").Value = -4999962
$ws.Range("H26").Value = 22500
$ws.Range("I26").Value = 5000
$ws.Range("K26").Value = 5000
$ws.Range("M26").Value = -4708
$ws.Range("H86").Value = 27805.105
$ws.Range("I86").Value = 36978.57
$ws.Range("J86").Value = 2119.4
$ws.Range("K86").Value = 36978.57
$ws.Range("L86").Value = 2119.4
$ws.Range("M86").Value = -35855.57
$ws.Range("N86").Value = -4365.4
$ws.Range("H89").Value = 27805.105
$ws.Range("I89").Value = 36978.57
$ws.Range("J89").Value = 2119.4
$ws.Range("K89").Value = 184892.85
$ws.Range("L89").Value = 10597
$ws.Range("M89").Value = -179276.85
$ws.Range("N89").Value = -21829
$ws.Range("H94").Value = 4069.353
$ws.Range("I94").Value = 1574.8334
$ws.Range("K94").Value = 1574.8334
$ws.Range("M94").Value = -1123.8334
$ws.Range("H105").Value = 1721.6111
$ws.Range("I105").Value = 1489.909
$ws.Range("K105").Value = 1489.909
$ws.Range("M105").Value = 257.0909999999999
$ws.Range("H134").Value = 2963.8572
$ws.Range("I134").Value = 3002.3076
$ws.Range("J134").Value = 2464
$ws.Range("K134").Value = 9006.9228
$ws.Range("L134").Value = 7392
$ws.Range("M134").Value = -6471.9228
$ws.Range("N134").Value = -12462
$ws.Range("H138").Value = 40162.223
$ws.Range("J138").Value = 40162.223
$ws.Range("L138").Value = 40162.223
$ws.Range("N138").Value = -50442.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 14407.892
$ws.Range("I58").Value = 1062.2693
$ws.Range("J58").Value = 45952.09
$ws.Range("K58").Value = 1062.2693
$ws.Range("L58").Value = 45952.09
$ws.Range("M58").Value = -859.2692999999999
$ws.Range("N58").Value = -46358.09
$ws.Range("H136").Value = 14407.892
$ws.Range("I136").Value = 1062.2693
$ws.Range("J136").Value = 45952.09
$ws.Range("K136").Value = 3186.8079
$ws.Range("L136").Value = 137856.27
$ws.Range("M136").Value = -636.8078999999998
$ws.Range("N136").Value = -142956.27

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 650.675
$ws.Range("I5").Value = 497.57895
$ws.Range("K5").Value = 1492.73685
$ws.Range("M5").Value = -1380.73685
$ws.Range("H25").Value = 375
$ws.Range("I25").Value = 375
$ws.Range("K25").Value = 1125
$ws.Range("M25").Value = -956
$ws.Range("H30").Value = 375
$ws.Range("I30").Value = 375
$ws.Range("K30").Value = 1125
$ws.Range("M30").Value = -1023
$ws.Range("H68").Value = 1170.8148
$ws.Range("I68").Value = 574.75
$ws.Range("J68").Value = 1421.7894
$ws.Range("K68").Value = 1724.25
$ws.Range("L68").Value = 4265.3682
$ws.Range("M68").Value = -913.25
$ws.Range("N68").Value = -5887.3682
$ws.Range("H71").Value = 1170.8148
$ws.Range("I71").Value = 574.75
$ws.Range("J71").Value = 1421.7894
$ws.Range("K71").Value = 5172.75
$ws.Range("L71").Value = 12796.1046
$ws.Range("M71").Value = -1116.75
$ws.Range("N71").Value = -20908.1046
$ws.Range("H107").Value = 4475.6924
$ws.Range("J107").Value = 772.46155
$ws.Range("L107").Value = 2317.38465
$ws.Range("N107").Value = -6157.38465
$ws.Range("H113").Value = 400.5
$ws.Range("I113").Value = 431.66666
$ws.Range("J113").Value = 381.8
$ws.Range("K113").Value = 1294.99998
$ws.Range("L113").Value = 1145.4
$ws.Range("M113").Value = 875.0000199999999
$ws.Range("N113").Value = -5485.4
$ws.Range("H131").Value = 107158.39
$ws.Range("I131").Value = 453
$ws.Range("J131").Value = 119861.414
$ws.Range("K131").Value = 1359
$ws.Range("L131").Value = 359584.242
$ws.Range("M131").Value = 3681
$ws.Range("N131").Value = -369664.242
$ws.Range("H135").Value = 650.675
$ws.Range("I135").Value = 497.57895
$ws.Range("K135").Value = 4478.21055
$ws.Range("M135").Value = -1943.21055

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9749.0625
$ws.Range("I80").Value = 19367.5
$ws.Range("J80").Value = 3978
$ws.Range("K80").Value = 19367.5
$ws.Range("L80").Value = 3978
$ws.Range("M80").Value = -18369.5
$ws.Range("N80").Value = -5974
$ws.Range("H83").Value = 9749.0625
$ws.Range("I83").Value = 19367.5
$ws.Range("J83").Value = 3978
$ws.Range("K83").Value = 96837.5
$ws.Range("L83").Value = 19890
$ws.Range("M83").Value = -91845.5
$ws.Range("N83").Value = -29874
$ws.Range("H113").Value = 3717.95
$ws.Range("I113").Value = 3017.2666
$ws.Range("J113").Value = 5820
$ws.Range("K113").Value = 3017.2666
$ws.Range("L113").Value = 5820
$ws.Range("M113").Value = -847.2665999999999
$ws.Range("N113").Value = -10160
$ws.Range("H122").Value = 4119.0625
$ws.Range("I122").Value = 3545.182
$ws.Range("K122").Value = 10635.546
$ws.Range("M122").Value = -8185.545999999998
$ws.Range("H126").Value = 5282.212
$ws.Range("I126").Value = 5259.95
$ws.Range("K126").Value = 15779.85
$ws.Range("M126").Value = -13309.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 29900
$ws.Range("J133").Value = 29900
$ws.Range("L133").Value = 29900
$ws.Range("N133").Value = -34960
$ws.Range("H136").Value = 24991.227
$ws.Range("I136").Value = 34626.867
$ws.Range("J136").Value = 4343.4287
$ws.Range("K136").Value = 103880.601
$ws.Range("L136").Value = 13030.2861
$ws.Range("M136").Value = -101330.601
$ws.Range("N136").Value = -18130.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1787.5
$ws.Range("I81").Value = 2200
$ws.Range("K81").Value = 4400
$ws.Range("M81").Value = -3339
$ws.Range("H84").Value = 1787.5
$ws.Range("I84").Value = 2200
$ws.Range("K84").Value = 22000
$ws.Range("M84").Value = -16696
$ws.Range("H94").Value = 25628.334
$ws.Range("J94").Value = 25628.334
$ws.Range("L94").Value = 25628.334
$ws.Range("N94").Value = -27430.334
$ws.Range("H136").Value = 1069.8889
$ws.Range("I136").Value = 589.8570999999999
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 1769.5713
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = 780.4287000000002
$ws.Range("N136").Value = -13350

Write-Host "Applied all cell updates"
